# Auto-generated edit script: apply scheduled market-price refresh values
# to the Leve profit calculation sheets (columns H-N) per the commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2689.2727
$ws.Range("I40").Value = 2300
$ws.Range("J40").Value = 2728.2
$ws.Range("K40").Value = 2300
$ws.Range("L40").Value = 2728.2
$ws.Range("M40").Value = -2125
$ws.Range("N40").Value = -3078.2
$ws.Range("H49").Value = 579.8
$ws.Range("I49").Value = 300
$ws.Range("J49").Value = 999.5
$ws.Range("K49").Value = 900
$ws.Range("L49").Value = 2998.5
$ws.Range("M49").Value = -764
$ws.Range("N49").Value = -3270.5
$ws.Range("H53").Value = 245.84616
$ws.Range("I53").Value = 226.125
$ws.Range("K53").Value = 226.125
$ws.Range("M53").Value = 410.875
$ws.Range("H93").Value = 35601
$ws.Range("J93").Value = 35601
$ws.Range("L93").Value = 35601
$ws.Range("N93").Value = -40593
$ws.Range("H112").Value = 1543.52
$ws.Range("J112").Value = 1543.52
$ws.Range("L112").Value = 4630.559999999999
$ws.Range("N112").Value = -6846.559999999999
$ws.Range("H129").Value = 2355.5
$ws.Range("I129").Value = 17224.5
$ws.Range("J129").Value = 916.5645
$ws.Range("K129").Value = 51673.5
$ws.Range("L129").Value = 2749.6935
$ws.Range("M129").Value = -46673.5
$ws.Range("N129").Value = -12749.6935
$ws.Range("H132").Value = 5562712.5
$ws.Range("I132").Value = 6104733
$ws.Range("J132").Value = 6999.75
$ws.Range("K132").Value = 18314199
$ws.Range("L132").Value = 20999.25
$ws.Range("M132").Value = -18311669
$ws.Range("N132").Value = -26059.25
$ws.Range("H137").Value = 1505.6897
$ws.Range("I137").Value = 1325.9131
$ws.Range("J137").Value = 2194.8333
$ws.Range("K137").Value = 3977.7393
$ws.Range("L137").Value = 6584.499899999999
$ws.Range("M137").Value = -1427.7393
$ws.Range("N137").Value = -11684.4999
$ws.Range("H138").Value = 2469.121
$ws.Range("I138").Value = 1109.9565
$ws.Range("J138").Value = 2880.4473
$ws.Range("K138").Value = 3329.8695
$ws.Range("L138").Value = 8641.341899999999
$ws.Range("M138").Value = 1810.1305
$ws.Range("N138").Value = -18921.3419

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29258.15
$ws.Range("I32").Value = 11697.912
$ws.Range("K32").Value = 11697.912
$ws.Range("M32").Value = -11410.912
$ws.Range("H95").Value = 25441.6
$ws.Range("J95").Value = 25441.6
$ws.Range("L95").Value = 25441.6
$ws.Range("N95").Value = -30933.6
$ws.Range("H132").Value = 36192.11
$ws.Range("I132").Value = 48573.848
$ws.Range("J132").Value = 3999.6
$ws.Range("K132").Value = 145721.544
$ws.Range("L132").Value = 11998.8
$ws.Range("M132").Value = -143191.544
$ws.Range("N132").Value = -17058.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 467.66666
$ws.Range("I7").Value = 467.66666
$ws.Range("K7").Value = 467.66666
$ws.Range("M7").Value = -354.66666
$ws.Range("H33").Value = 4264.2
$ws.Range("I33").Value = 2840.3333
$ws.Range("J33").Value = 6400
$ws.Range("K33").Value = 2840.3333
$ws.Range("L33").Value = 6400
$ws.Range("M33").Value = -2504.3333
$ws.Range("N33").Value = -7072
$ws.Range("H44").Value = 11428.571
$ws.Range("I44").Value = 10000
$ws.Range("K44").Value = 10000
$ws.Range("M44").Value = -9503
$ws.Range("H54").Value = 6694.8
$ws.Range("I54").Value = 4321.143
$ws.Range("J54").Value = 12233.333
$ws.Range("K54").Value = 4321.143
$ws.Range("L54").Value = 12233.333
$ws.Range("M54").Value = -3837.143
$ws.Range("N54").Value = -13201.333
$ws.Range("H107").Value = 125058240
$ws.Range("I107").Value = 200092270
$ws.Range("J107").Value = 1533.3334
$ws.Range("K107").Value = 200092270
$ws.Range("L107").Value = 1533.3334
$ws.Range("M107").Value = -200090350
$ws.Range("N107").Value = -5373.3334
$ws.Range("H134").Value = 3274.7112
$ws.Range("I134").Value = 3484.6333
$ws.Range("J134").Value = 2854.8667
$ws.Range("K134").Value = 10453.8999
$ws.Range("L134").Value = 8564.6001
$ws.Range("M134").Value = -7918.8999
$ws.Range("N134").Value = -13634.6001
$ws.Range("H141").Value = 54920
$ws.Range("J141").Value = 59900
$ws.Range("L141").Value = 59900
$ws.Range("N141").Value = -70260

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 143929
$ws.Range("I16").Value = 1300
$ws.Range("K16").Value = 1300
$ws.Range("M16").Value = -1013
$ws.Range("H31").Value = 57595.19
$ws.Range("I31").Value = 1347.3334
$ws.Range("J31").Value = 105807.64
$ws.Range("K31").Value = 1347.3334
$ws.Range("L31").Value = 105807.64
$ws.Range("M31").Value = -1052.3334
$ws.Range("N31").Value = -106397.64
$ws.Range("H34").Value = 57595.19
$ws.Range("I34").Value = 1347.3334
$ws.Range("J34").Value = 105807.64
$ws.Range("K34").Value = 1347.3334
$ws.Range("L34").Value = 105807.64
$ws.Range("M34").Value = -1145.3334
$ws.Range("N34").Value = -106211.64
$ws.Range("H113").Value = 143929
$ws.Range("I113").Value = 1300
$ws.Range("K113").Value = 1300
$ws.Range("M113").Value = 870
$ws.Range("H134").Value = 1209.174
$ws.Range("I134").Value = 678.2222
$ws.Range("J134").Value = 3120.6
$ws.Range("K134").Value = 2034.6666
$ws.Range("L134").Value = 9361.799999999999
$ws.Range("M134").Value = 500.3334
$ws.Range("N134").Value = -14431.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 40.347828
$ws.Range("J12").Value = 48.11111
$ws.Range("L12").Value = 144.33333
$ws.Range("N12").Value = -490.33333
$ws.Range("H37").Value = 522805.9
$ws.Range("J37").Value = 522805.9
$ws.Range("L37").Value = 1568417.7
$ws.Range("N37").Value = -1568641.7
$ws.Range("H86").Value = 1000
$ws.Range("I86").Value = 1233.3334
$ws.Range("J86").Value = 825
$ws.Range("K86").Value = 3700.0002
$ws.Range("L86").Value = 2475
$ws.Range("M86").Value = -2514.0002
$ws.Range("N86").Value = -4847
$ws.Range("H89").Value = 1000
$ws.Range("I89").Value = 1233.3334
$ws.Range("J89").Value = 825
$ws.Range("K89").Value = 11100.0006
$ws.Range("L89").Value = 7425
$ws.Range("M89").Value = -5172.000599999999
$ws.Range("N89").Value = -19281
$ws.Range("H92").Value = 1194.3334
$ws.Range("I92").Value = 290
$ws.Range("J92").Value = 3003
$ws.Range("K92").Value = 870
$ws.Range("L92").Value = 9009
$ws.Range("M92").Value = 378
$ws.Range("N92").Value = -11505
$ws.Range("H98").Value = 101057.6
$ws.Range("I98").Value = 3
$ws.Range("J98").Value = 112285.89
$ws.Range("K98").Value = 9
$ws.Range("L98").Value = 336857.67
$ws.Range("M98").Value = 1489
$ws.Range("N98").Value = -339853.67
$ws.Range("H113").Value = 1024.5
$ws.Range("I113").Value = 1833.5
$ws.Range("J113").Value = 620
$ws.Range("K113").Value = 5500.5
$ws.Range("L113").Value = 1860
$ws.Range("M113").Value = -3330.5
$ws.Range("N113").Value = -6200
$ws.Range("H122").Value = 697.5909
$ws.Range("I122").Value = 546.3333
$ws.Range("J122").Value = 754.3125
$ws.Range("K122").Value = 4916.9997
$ws.Range("L122").Value = 6788.8125
$ws.Range("M122").Value = -2466.9997
$ws.Range("N122").Value = -11688.8125
$ws.Range("H131").Value = 839.22
$ws.Range("I131").Value = 581.6
$ws.Range("J131").Value = 884.6824
$ws.Range("K131").Value = 1744.8
$ws.Range("L131").Value = 2654.0472
$ws.Range("M131").Value = 3295.2
$ws.Range("N131").Value = -12734.0472
$ws.Range("H132").Value = 371938.25
$ws.Range("I132").Value = 811.4
$ws.Range("J132").Value = 590248.2
$ws.Range("K132").Value = 7302.599999999999
$ws.Range("L132").Value = 5312233.8
$ws.Range("M132").Value = -4772.599999999999
$ws.Range("N132").Value = -5317293.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H132").Value = 3162.75
$ws.Range("I132").Value = 2335.375
$ws.Range("J132").Value = 4265.9165
$ws.Range("K132").Value = 7006.125
$ws.Range("L132").Value = 12797.7495
$ws.Range("M132").Value = -4476.125
$ws.Range("N132").Value = -17857.7495
$ws.Range("H134").Value = 15661.5
$ws.Range("J134").Value = 15661.5
$ws.Range("L134").Value = 46984.5
$ws.Range("N134").Value = -52054.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4257.696
$ws.Range("I132").Value = 5063.769
$ws.Range("J132").Value = 3209.8
$ws.Range("K132").Value = 15191.307
$ws.Range("L132").Value = 9629.400000000001
$ws.Range("M132").Value = -12661.307
$ws.Range("N132").Value = -14689.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 1200
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1200
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 1200
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -1754
$ws.Range("H81").Value = 225808.78
$ws.Range("I81").Value = 1000000
$ws.Range("J81").Value = 129034.875
$ws.Range("K81").Value = 2000000
$ws.Range("L81").Value = 258069.75
$ws.Range("M81").Value = -1998939
$ws.Range("N81").Value = -260191.75
$ws.Range("H84").Value = 225808.78
$ws.Range("I84").Value = 1000000
$ws.Range("J84").Value = 129034.875
$ws.Range("K84").Value = 10000000
$ws.Range("L84").Value = 1290348.75
$ws.Range("M84").Value = -9994696
$ws.Range("N84").Value = -1300956.75
$ws.Range("H136").Value = 22168.445
$ws.Range("I136").Value = 50714.8
$ws.Range("J136").Value = 5376.4707
$ws.Range("K136").Value = 152144.4
$ws.Range("L136").Value = 16129.4121
$ws.Range("M136").Value = -149594.4
$ws.Range("N136").Value = -21229.4121
